# Auto-generated edit script applying cryptos.xlsx price/volume/ranking updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.954.18"
$ws.Range("E2").Value = "'  -0.49%  "
$ws.Range("D3").Value = "'2.043.82"
$ws.Range("E3").Value = "'  -0.14%  "
$ws.Range("E4").Value = "'  -0.25%  "
$ws.Range("D5").Value = "'248.07"
$ws.Range("E5").Value = "'  -0.47%  "
$ws.Range("E6").Value = "'  +0.39%  "
$ws.Range("E7").Value = "'  +0.07%  "
$ws.Range("D8").Value = "'55.75"
$ws.Range("E8").Value = "'  -1.21%  "
$ws.Range("E9").Value = "'  -0.75%  "
$ws.Range("E10").Value = "'  +2.14%  "
$ws.Range("E11").Value = "'  +1.46%  "
$ws.Range("D12").Value = "'15.71"
$ws.Range("E12").Value = "'  +3.48%  "
$ws.Range("D13").Value = "'2.339.61"
$ws.Range("E13").Value = "'  -0.11%  "
$ws.Range("D14").Value = "'5.63"
$ws.Range("E14").Value = "'  +6.20%  "
$ws.Range("D15").Value = "'0.787"
$ws.Range("E15").Value = "'  -5.06%  "
$ws.Range("D16").Value = "'2.040.87"
$ws.Range("E16").Value = "'  -0.21%  "
$ws.Range("D17").Value = "'36.932.50"
$ws.Range("D18").Value = "'16.39"
$ws.Range("E18").Value = "'  +12.88%  "
$ws.Range("D19").Value = "'73.59"
$ws.Range("E19").Value = "'  +1.25%  "
$ws.Range("D20").Value = "'0.0₃0892"
$ws.Range("E20").Value = "'  +2.50%  "
$ws.Range("D21").Value = "'5.30"
$ws.Range("E21").Value = "'  +0.51%  "
$ws.Range("D22").Value = "'235.27"
$ws.Range("E22").Value = "'  -1.39%  "
$ws.Range("E23").Value = "'  -0.03%  "
$ws.Range("D24").Value = "'2.35"
$ws.Range("E24").Value = "'  -3.16%  "
$ws.Range("D25").Value = "'2.18"
$ws.Range("E25").Value = "'  +8.31%  "
$ws.Range("D26").Value = "'167.27"
$ws.Range("E26").Value = "'  -1.95%  "
$ws.Range("D27").Value = "'9.05"
$ws.Range("E27").Value = "'  -1.57%  "
$ws.Range("D28").Value = "'19.63"
$ws.Range("E28").Value = "'  -3.94%  "
$ws.Range("E29").Value = "'  +0.35%  "
$ws.Range("E30").Value = "'  +1.65%  "
$ws.Range("E31").Value = "'  +1.03%  "
$ws.Range("D32").Value = "'0.0608"
$ws.Range("E32").Value = "'  -4.09%  "
$ws.Range("D33").Value = "'4.39"
$ws.Range("E33").Value = "'  +0.10%  "
$ws.Range("E34").Value = "'  -0.07%  "
$ws.Range("D35").Value = "'0.0870"
$ws.Range("E35").Value = "'  +1.46%  "
$ws.Range("E36").Value = "'  -4.71%  "
$ws.Range("E37").Value = "'  -1.74%  "
$ws.Range("B38").Value = "'Cronos"
$ws.Range("C38").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D38").Value = "'0.106"
$ws.Range("E38").Value = "'  -1.92%  "
$ws.Range("B39").Value = "'HuobiToken"
$ws.Range("C39").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D39").Value = "'3.22"
$ws.Range("E39").Value = "'  +15.39%  "
$ws.Range("B40").Value = "'TrustWalletToken"
$ws.Range("C40").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.33"
$ws.Range("E40").Value = "'  -1.85%  "
$ws.Range("B41").Value = "'THORChain"
$ws.Range("C41").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").Value = "'4.86"
$ws.Range("E41").Value = "'  +20.85%  "
$ws.Range("B42").Value = "'VeChain"
$ws.Range("C42").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0219"
$ws.Range("E42").Value = "'  -3.17%  "
$ws.Range("D43").Value = "'17.16"
$ws.Range("E43").Value = "'  -5.82%  "
$ws.Range("D44").Value = "'95.02"
$ws.Range("E44").Value = "'  -2.76%  "
$ws.Range("E45").Value = "'  -3.67%  "
$ws.Range("E46").Value = "'  +0.91%  "
$ws.Range("D47").Value = "'1.272.46"
$ws.Range("E47").Value = "'  -2.67%  "
$ws.Range("E48").Value = "'  -2.28%  "
$ws.Range("D49").Value = "'2.225.78"
$ws.Range("E50").Value = "'  -3.70%  "
$ws.Range("D51").Value = "'41.83"
$ws.Range("E51").Value = "'  -7.97%  "
